$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "ACYDBNjkP-skmYhAC1ElMPOvzbtqF-UM8KlAjRXaL5jfkFjX23TF88gMaMNLvh5m09c79Ys"
$ws.Range("C2").Value = "2024-11-02T01:20:11.877Z"
$ws.Range("D2").Value = "2024-11-02T01:20:11.877676Z"
$ws.Range("E2").Value = "'10"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = "'3"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = "adsf"
$ws.Range("H2").Value = "Very good"
$ws.Range("I2").Value = "Yes"
$ws.Range("J2").Value = "Very good"
$ws.Range("K2").Value = "Very good"
$ws.Range("L2").Value = "asdf"
$ws.Range("M2").Value = "Kevin.LopezChavez01@student.csulb.edu"
$ws.Range("N2").Value = "'7"
$ws.Range("N2").Style = "Normal"
$ws.Range("O2").Value = "asdf"
$ws.Range("P2").Value = "Satisfactory"
$ws.Range("Q2").Value = "Very good"
$ws.Range("R2").Value = "Satisfactory"
$ws.Range("S2").Value = "Yes"
$ws.Range("T2").Value = "Very good"
